$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Duplicate the block of rows 59-74 down to rows 79-94 -----------------
$ws.Range("A59:B74").Copy($ws.Range("A79:B94"))

# --- 2) Row 30 ("How many scales can I define?" / "...define 23 scales.")
#        gets duplicated (twice) into rows 95-96 ------------------------------
$ws.Range("A30:B30").Copy($ws.Range("A95:B95"))
$ws.Range("A30:B30").Copy($ws.Range("A96:B96"))

# --- 3) Duplicate rows 77-78 down to rows 97-98 -------------------------------
$ws.Range("A77:B78").Copy($ws.Range("A97:B98"))

# --- 4) Append three brand-new question/answer pairs (each appearing twice) --
$newQA = @(
    @{ Q = "What is 10 * 10?"; A = "The answer is 100." },
    @{ Q = "How to plot a neutron density log?"; A = "To plot a neutron density log in the GEO application, follow these steps:`n1. Open the well log data in the GEO application.`n2. Select the neutron density log data and click on the `"Plot`" button.`n3. Choose the desired plot type (e.g., log-log, semi-log) and adjust the settings as needed.`n4. Click `"OK`" to generate the plot.`nNote: The exact steps may vary depending on the specific version of the GEO application being used." },
    @{ Q = "Summarise the types of limits in the GEO application."; A = "The GEO application supports two types of limits:`n1. Hard limits`n2. Soft limits" }
)

# Write every new question first (so they become consecutive new shared-string
# entries), then every new answer (also consecutive), matching the way the
# question/answer string pools are appended to in the workbook's source data.
$qRows = @(99, 101, 103)
for ($i = 0; $i -lt $newQA.Count; $i++) {
    $r = $qRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $newQA[$i].Q
    $ws.Cells.Item($r + 1, 1).Value2 = $newQA[$i].Q
}
for ($i = 0; $i -lt $newQA.Count; $i++) {
    $r = $qRows[$i]
    $ws.Cells.Item($r, 2).Value2 = $newQA[$i].A
    $ws.Cells.Item($r + 1, 2).Value2 = $newQA[$i].A
}
